# Update NATMI ligand-receptor edge-weight table with freshly recomputed
# TPM-based statistics (ligand/receptor expressing-cell counts, detection
# rates, average/total expression, derived specificities and edge weights).
# Only the data cells change; headers, labels and formatting are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> hashtable of column letter -> new value
$updates = [ordered]@{
    2  = [ordered]@{
        E = 3; F = 1; G = 1.644726333333334; H = 4.934179
        I = 0.03084360558270512; J = 0.03084360558270512
        K = 3; L = 1; M = 11.608856; N = 34.826568
        O = 0.6986105367350092; P = 0.6986105367350093
        Q = 19.09339116307467; R = 171.840520467672
        S = 0.02154766785097655; T = 0.02154766785097655
    }
    3  = [ordered]@{
        E = 3; F = 1; G = 1.644726333333334; H = 4.934179
        I = 0.03084360558270512; J = 0.03084360558270512
        O = 0.204840297499807; P = 0.204840297499807
        Q = 5.598392409600779; R = 50.385531686407
        S = 0.006318013343528025; T = 0.006318013343528026
    }
    4  = [ordered]@{
        E = 3; F = 1; G = 1.644726333333334; H = 4.934179
        I = 0.03084360558270512; J = 0.03084360558270512
        O = 0.0965491657651837; P = 0.09654916576518371
        Q = 2.63873917080989; R = 23.748652537289
        S = 0.002977924388200542; T = 0.002977924388200542
    }
    5  = [ordered]@{
        I = 0.828024694817689; J = 0.828024694817689
        K = 3; L = 1; M = 11.608856; N = 34.826568
        O = 0.6986105367350092; P = 0.6986105367350093
        Q = 512.5794825915121; R = 4613.215343323609
        S = 0.5784667764764279; T = 0.578466776476428
    }
    6  = [ordered]@{
        I = 0.828024694817689; J = 0.828024694817689
        O = 0.204840297499807; P = 0.204840297499807
        S = 0.1696128248236423; T = 0.1696128248236423
    }
    7  = [ordered]@{
        I = 0.828024694817689; J = 0.828024694817689
        O = 0.0965491657651837; P = 0.09654916576518371
        S = 0.07994509351761871; T = 0.07994509351761872
    }
    8  = [ordered]@{
        I = 0.1411316995996059; J = 0.1411316995996059
        K = 3; L = 1; M = 11.608856; N = 34.826568
        O = 0.6986105367350092; P = 0.6986105367350093
        Q = 87.36600974679202; R = 786.2940877211281
        S = 0.09859609240760475; T = 0.09859609240760477
    }
    9  = [ordered]@{
        I = 0.1411316995996059; J = 0.1411316995996059
        O = 0.204840297499807; P = 0.204840297499807
        S = 0.02890945933263666; T = 0.02890945933263667
    }
    10 = [ordered]@{
        I = 0.1411316995996059; J = 0.1411316995996059
        O = 0.0965491657651837; P = 0.09654916576518371
        S = 0.01362614785936446; T = 0.01362614785936446
    }
}

foreach ($row in $updates.Keys) {
    $rowUpdates = $updates[$row]
    foreach ($col in $rowUpdates.Keys) {
        $ws.Range("$col$row").Value = $rowUpdates[$col]
    }
}
